$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Finnish Veikkausliiga | Ilves vs KuPS
$ws.Cells.Item(2,1).Value = 'Finnish Veikkausliiga'
$ws.Cells.Item(2,2).NumberFormat = '@'
$ws.Cells.Item(2,2).Value = '2025-10-16'
$ws.Cells.Item(2,2).Style = 'Normal'
$ws.Cells.Item(2,3).Value = '12:00:00'
$ws.Cells.Item(2,4).Value = 'Ilves'
$ws.Cells.Item(2,5).Value = 'KuPS'
$ws.Cells.Item(2,6).Value = 2.1
$ws.Cells.Item(2,7).Value = 2.14
$ws.Cells.Item(2,8).Value = 3.5
$ws.Cells.Item(2,9).Value = 3.6
$ws.Cells.Item(2,10).Value = 3.95
$ws.Cells.Item(2,11).Value = 4.1
$ws.Cells.Item(2,12).Value = 1.01
$ws.Cells.Item(2,13).Value = 1.03
$ws.Cells.Item(2,14).Value = 6.2
$ws.Cells.Item(2,15).Value = 1.16
$ws.Cells.Item(2,16).Value = 2.76
$ws.Cells.Item(2,17).Value = 1.51
$ws.Cells.Item(2,18).Value = 1.72
$ws.Cells.Item(2,19).Value = 2.28
$ws.Cells.Item(2,20).Value = 1.48
$ws.Cells.Item(2,21).Value = 2.8
$ws.Cells.Item(2,22).Value = 1.38
$ws.Cells.Item(2,23).Value = 1.87
$ws.Cells.Item(2,24).Value = 28
$ws.Cells.Item(2,25).Value = 980
$ws.Cells.Item(2,26).Value = 32
$ws.Cells.Item(2,27).Value = 980
$ws.Cells.Item(2,28).Value = 19.5
$ws.Cells.Item(2,29).Value = 11
$ws.Cells.Item(2,30).Value = 15.5
$ws.Cells.Item(2,31).Value = 34
$ws.Cells.Item(2,32).Value = 20
$ws.Cells.Item(2,33).Value = 11.5
$ws.Cells.Item(2,34).Value = 15.5
$ws.Cells.Item(2,35).Value = 36
$ws.Cells.Item(2,36).Value = 29
$ws.Cells.Item(2,37).Value = 18.5
$ws.Cells.Item(2,38).Value = 26
$ws.Cells.Item(2,39).Value = 1000
$ws.Cells.Item(2,40).Value = 10
$ws.Cells.Item(2,41).Value = 30

# Row 3: Danish 1st Division | Esbjerg vs Hobro
$ws.Cells.Item(3,1).Value = 'Danish 1st Division'
$ws.Cells.Item(3,2).NumberFormat = '@'
$ws.Cells.Item(3,2).Value = '2025-10-16'
$ws.Cells.Item(3,2).Style = 'Normal'
$ws.Cells.Item(3,3).Value = '13:30:00'
$ws.Cells.Item(3,4).Value = 'Esbjerg'
$ws.Cells.Item(3,5).Value = 'Hobro'
$ws.Cells.Item(3,6).Value = 1.97
$ws.Cells.Item(3,7).Value = 980
$ws.Cells.Item(3,8).Value = 1.32
$ws.Cells.Item(3,9).Value = 980
$ws.Cells.Item(3,10).Value = 1.36
$ws.Cells.Item(3,11).Value = 950
$ws.Cells.Item(3,12).Value = 1.01
$ws.Cells.Item(3,13).Value = 1.04
$ws.Cells.Item(3,14).Value = 1.78
$ws.Cells.Item(3,15).Value = 1.21
$ws.Cells.Item(3,16).Value = 1.78
$ws.Cells.Item(3,17).Value = 1.21
$ws.Cells.Item(3,18).Value = 1.45
$ws.Cells.Item(3,19).Value = 1.21
$ws.Cells.Item(3,20).Value = 1.03
$ws.Cells.Item(3,21).Value = 1.03
$ws.Cells.Item(3,22).Value = 1.34
$ws.Cells.Item(3,23).Value = 1.01
$ws.Cells.Item(3,24).Value = 1000
$ws.Cells.Item(3,25).Value = 980
$ws.Cells.Item(3,26).Value = 980
$ws.Cells.Item(3,27).Value = 1000
$ws.Cells.Item(3,28).Value = 980
$ws.Cells.Item(3,29).Value = 980
$ws.Cells.Item(3,30).Value = 980
$ws.Cells.Item(3,31).Value = 980
$ws.Cells.Item(3,32).Value = 980
$ws.Cells.Item(3,33).Value = 12.5
$ws.Cells.Item(3,34).Value = 980
$ws.Cells.Item(3,35).Value = 1000
$ws.Cells.Item(3,36).Value = 980
$ws.Cells.Item(3,37).Value = 980
$ws.Cells.Item(3,38).Value = 980
$ws.Cells.Item(3,39).Value = 1000
$ws.Cells.Item(3,40).Value = 1000
$ws.Cells.Item(3,41).Value = 1000

# Row 4: Romanian Liga I | Csikszereda vs CFR Cluj
$ws.Cells.Item(4,1).Value = 'Romanian Liga I'
$ws.Cells.Item(4,2).NumberFormat = '@'
$ws.Cells.Item(4,2).Value = '2025-10-16'
$ws.Cells.Item(4,2).Style = 'Normal'
$ws.Cells.Item(4,3).Value = '14:30:00'
$ws.Cells.Item(4,4).Value = 'Csikszereda'
$ws.Cells.Item(4,5).Value = 'CFR Cluj'
$ws.Cells.Item(4,6).Value = 2.14
$ws.Cells.Item(4,7).Value = 1000
$ws.Cells.Item(4,8).Value = 1.64
$ws.Cells.Item(4,9).Value = 1.76
$ws.Cells.Item(4,10).Value = 4
$ws.Cells.Item(4,11).Value = 4.3
$ws.Cells.Item(4,12).Value = 1.01
$ws.Cells.Item(4,13).Value = 1.05
$ws.Cells.Item(4,14).Value = 1.11
$ws.Cells.Item(4,15).Value = 1.29
$ws.Cells.Item(4,16).Value = 1.38
$ws.Cells.Item(4,17).Value = 1.85
$ws.Cells.Item(4,18).Value = 1.18
$ws.Cells.Item(4,19).Value = 1.05
$ws.Cells.Item(4,20).Value = 1.03
$ws.Cells.Item(4,21).Value = 2
$ws.Cells.Item(4,22).Value = 2.3
$ws.Cells.Item(4,23).Value = 1.18
$ws.Cells.Item(4,24).Value = 1000
$ws.Cells.Item(4,25).Value = 11.5
$ws.Cells.Item(4,26).Value = 13.5
$ws.Cells.Item(4,27).Value = 22
$ws.Cells.Item(4,28).Value = 27
$ws.Cells.Item(4,29).Value = 980
$ws.Cells.Item(4,30).Value = 13
$ws.Cells.Item(4,31).Value = 28
$ws.Cells.Item(4,32).Value = 1000
$ws.Cells.Item(4,33).Value = 30
$ws.Cells.Item(4,34).Value = 980
$ws.Cells.Item(4,35).Value = 980
$ws.Cells.Item(4,36).Value = 1000
$ws.Cells.Item(4,37).Value = 1000
$ws.Cells.Item(4,38).Value = 1000
$ws.Cells.Item(4,39).Value = 1000
$ws.Cells.Item(4,40).Value = 1000
$ws.Cells.Item(4,41).Value = 13.5

# Row 5: Paraguayan Primera Division | General Caballero vs Club Atletico Tembetary
$ws.Cells.Item(5,1).Value = 'Paraguayan Primera Division'
$ws.Cells.Item(5,2).NumberFormat = '@'
$ws.Cells.Item(5,2).Value = '2025-10-16'
$ws.Cells.Item(5,2).Style = 'Normal'
$ws.Cells.Item(5,3).Value = '18:30:00'
$ws.Cells.Item(5,4).Value = 'General Caballero'
$ws.Cells.Item(5,5).Value = 'Club Atletico Tembetary'
$ws.Cells.Item(5,6).Value = 0
$ws.Cells.Item(5,7).Value = 0
$ws.Cells.Item(5,8).Value = 0
$ws.Cells.Item(5,9).Value = 0
$ws.Cells.Item(5,10).Value = 0
$ws.Cells.Item(5,11).Value = 0
$ws.Cells.Item(5,12).Value = 1.01
$ws.Cells.Item(5,13).Value = 1.01
$ws.Cells.Item(5,14).Value = 1.11
$ws.Cells.Item(5,15).Value = 1.4
$ws.Cells.Item(5,16).Value = 1.24
$ws.Cells.Item(5,17).Value = 1.01
$ws.Cells.Item(5,18).Value = 1.18
$ws.Cells.Item(5,19).Value = 1.05
$ws.Cells.Item(5,20).Value = 1.03
$ws.Cells.Item(5,21).Value = 1.03
$ws.Cells.Item(5,22).Value = 1.01
$ws.Cells.Item(5,23).Value = 1.01
$ws.Cells.Item(5,24).Value = 1000
$ws.Cells.Item(5,25).Value = 1000
$ws.Cells.Item(5,26).Value = 1000
$ws.Cells.Item(5,27).Value = 1000
$ws.Cells.Item(5,28).Value = 1000
$ws.Cells.Item(5,29).Value = 1000
$ws.Cells.Item(5,30).Value = 1000
$ws.Cells.Item(5,31).Value = 1000
$ws.Cells.Item(5,32).Value = 1000
$ws.Cells.Item(5,33).Value = 1000
$ws.Cells.Item(5,34).Value = 1000
$ws.Cells.Item(5,35).Value = 1000
$ws.Cells.Item(5,36).Value = 1000
$ws.Cells.Item(5,37).Value = 1000
$ws.Cells.Item(5,38).Value = 1000
$ws.Cells.Item(5,39).Value = 1000
$ws.Cells.Item(5,40).Value = 1000
$ws.Cells.Item(5,41).Value = 1000

# Row 6: Brazilian Serie A | Gremio vs Sao Paulo
$ws.Cells.Item(6,1).Value = 'Brazilian Serie A'
$ws.Cells.Item(6,2).NumberFormat = '@'
$ws.Cells.Item(6,2).Value = '2025-10-16'
$ws.Cells.Item(6,2).Style = 'Normal'
$ws.Cells.Item(6,3).Value = '19:00:00'
$ws.Cells.Item(6,4).Value = 'Gremio'
$ws.Cells.Item(6,5).Value = 'Sao Paulo'
$ws.Cells.Item(6,6).Value = 2.8
$ws.Cells.Item(6,7).Value = 3.2
$ws.Cells.Item(6,8).Value = 2.66
$ws.Cells.Item(6,9).Value = 2.84
$ws.Cells.Item(6,10).Value = 3.15
$ws.Cells.Item(6,11).Value = 3.3
$ws.Cells.Item(6,12).Value = 1.59
$ws.Cells.Item(6,13).Value = 1.08
$ws.Cells.Item(6,14).Value = 2.36
$ws.Cells.Item(6,15).Value = 1.55
$ws.Cells.Item(6,16).Value = 1.25
$ws.Cells.Item(6,17).Value = 1.55
$ws.Cells.Item(6,18).Value = 1.18
$ws.Cells.Item(6,19).Value = 4.8
$ws.Cells.Item(6,20).Value = 1.03
$ws.Cells.Item(6,21).Value = 1.03
$ws.Cells.Item(6,22).Value = 1.54
$ws.Cells.Item(6,23).Value = 1.46
$ws.Cells.Item(6,24).Value = 1000
$ws.Cells.Item(6,25).Value = 1000
$ws.Cells.Item(6,26).Value = 1000
$ws.Cells.Item(6,27).Value = 1000
$ws.Cells.Item(6,28).Value = 1000
$ws.Cells.Item(6,29).Value = 1000
$ws.Cells.Item(6,30).Value = 1000
$ws.Cells.Item(6,31).Value = 1000
$ws.Cells.Item(6,32).Value = 1000
$ws.Cells.Item(6,33).Value = 1000
$ws.Cells.Item(6,34).Value = 1000
$ws.Cells.Item(6,35).Value = 1000
$ws.Cells.Item(6,36).Value = 1000
$ws.Cells.Item(6,37).Value = 1000
$ws.Cells.Item(6,38).Value = 1000
$ws.Cells.Item(6,39).Value = 1000
$ws.Cells.Item(6,40).Value = 1000
$ws.Cells.Item(6,41).Value = 1000

# Row 7: Honduras Liga Nacional | Olancho vs Juticalpa
$ws.Cells.Item(7,1).Value = 'Honduras Liga Nacional'
$ws.Cells.Item(7,2).NumberFormat = '@'
$ws.Cells.Item(7,2).Value = '2025-10-16'
$ws.Cells.Item(7,2).Style = 'Normal'
$ws.Cells.Item(7,3).Value = '20:15:00'
$ws.Cells.Item(7,4).Value = 'Olancho'
$ws.Cells.Item(7,5).Value = 'Juticalpa'
$ws.Cells.Item(7,6).Value = 1.04
$ws.Cells.Item(7,7).Value = 1000
$ws.Cells.Item(7,8).Value = 1.04
$ws.Cells.Item(7,9).Value = 1000
$ws.Cells.Item(7,10).Value = 1.01
$ws.Cells.Item(7,11).Value = 1000
$ws.Cells.Item(7,12).Value = 1.01
$ws.Cells.Item(7,13).Value = 1.01
$ws.Cells.Item(7,14).Value = 1.24
$ws.Cells.Item(7,15).Value = 1.25
$ws.Cells.Item(7,16).Value = 1.24
$ws.Cells.Item(7,17).Value = 1.25
$ws.Cells.Item(7,18).Value = 1.18
$ws.Cells.Item(7,19).Value = 1.25
$ws.Cells.Item(7,20).Value = 1.01
$ws.Cells.Item(7,21).Value = 1.01
$ws.Cells.Item(7,22).Value = 1.01
$ws.Cells.Item(7,23).Value = 1.01
$ws.Cells.Item(7,24).Value = 1000
$ws.Cells.Item(7,25).Value = 1000
$ws.Cells.Item(7,26).Value = 1000
$ws.Cells.Item(7,27).Value = 1000
$ws.Cells.Item(7,28).Value = 1000
$ws.Cells.Item(7,29).Value = 1000
$ws.Cells.Item(7,30).Value = 1000
$ws.Cells.Item(7,31).Value = 1000
$ws.Cells.Item(7,32).Value = 1000
$ws.Cells.Item(7,33).Value = 1000
$ws.Cells.Item(7,34).Value = 1000
$ws.Cells.Item(7,35).Value = 1000
$ws.Cells.Item(7,36).Value = 1000
$ws.Cells.Item(7,37).Value = 1000
$ws.Cells.Item(7,38).Value = 1000
$ws.Cells.Item(7,39).Value = 1000
$ws.Cells.Item(7,40).Value = 1000
$ws.Cells.Item(7,41).Value = 1000

# Row 8: Colombian Primera B | Quindio vs Tigres FC Zipaquira
$ws.Cells.Item(8,1).Value = 'Colombian Primera B'
$ws.Cells.Item(8,2).NumberFormat = '@'
$ws.Cells.Item(8,2).Value = '2025-10-16'
$ws.Cells.Item(8,2).Style = 'Normal'
$ws.Cells.Item(8,3).Value = '21:00:00'
$ws.Cells.Item(8,4).Value = 'Quindio'
$ws.Cells.Item(8,5).Value = 'Tigres FC Zipaquira'
$ws.Cells.Item(8,6).Value = 1.43
$ws.Cells.Item(8,7).Value = 1.61
$ws.Cells.Item(8,8).Value = 6.8
$ws.Cells.Item(8,9).Value = 1000
$ws.Cells.Item(8,10).Value = 3.9
$ws.Cells.Item(8,11).Value = 7.2
$ws.Cells.Item(8,12).Value = 1.01
$ws.Cells.Item(8,13).Value = 1.01
$ws.Cells.Item(8,14).Value = 1.25
$ws.Cells.Item(8,15).Value = 1.01
$ws.Cells.Item(8,16).Value = 1.25
$ws.Cells.Item(8,17).Value = 1.99
$ws.Cells.Item(8,18).Value = 1.13
$ws.Cells.Item(8,19).Value = 1.99
$ws.Cells.Item(8,20).Value = 1.03
$ws.Cells.Item(8,21).Value = 1.03
$ws.Cells.Item(8,22).Value = 1.01
$ws.Cells.Item(8,23).Value = 2.78
$ws.Cells.Item(8,24).Value = 17
$ws.Cells.Item(8,25).Value = 55
$ws.Cells.Item(8,26).Value = 870
$ws.Cells.Item(8,27).Value = 1000
$ws.Cells.Item(8,28).Value = 9.4
$ws.Cells.Item(8,29).Value = 28
$ws.Cells.Item(8,30).Value = 500
$ws.Cells.Item(8,31).Value = 1000
$ws.Cells.Item(8,32).Value = 10
$ws.Cells.Item(8,33).Value = 38
$ws.Cells.Item(8,34).Value = 560
$ws.Cells.Item(8,35).Value = 1000
$ws.Cells.Item(8,36).Value = 21
$ws.Cells.Item(8,37).Value = 85
$ws.Cells.Item(8,38).Value = 990
$ws.Cells.Item(8,39).Value = 1000
$ws.Cells.Item(8,40).Value = 40
$ws.Cells.Item(8,41).Value = 1000

# Row 9: Brazilian Serie A | EC Vitoria Salvador vs Bahia
$ws.Cells.Item(9,1).Value = 'Brazilian Serie A'
$ws.Cells.Item(9,2).NumberFormat = '@'
$ws.Cells.Item(9,2).Value = '2025-10-16'
$ws.Cells.Item(9,2).Style = 'Normal'
$ws.Cells.Item(9,3).Value = '21:30:00'
$ws.Cells.Item(9,4).Value = 'EC Vitoria Salvador'
$ws.Cells.Item(9,5).Value = 'Bahia'
$ws.Cells.Item(9,6).Value = 2.58
$ws.Cells.Item(9,7).Value = 3.5
$ws.Cells.Item(9,8).Value = 2.52
$ws.Cells.Item(9,9).Value = 870
$ws.Cells.Item(9,10).Value = 3.15
$ws.Cells.Item(9,11).Value = 3.2
$ws.Cells.Item(9,12).Value = 1.56
$ws.Cells.Item(9,13).Value = 1.07
$ws.Cells.Item(9,14).Value = 1.06
$ws.Cells.Item(9,15).Value = 1.5
$ws.Cells.Item(9,16).Value = 1.25
$ws.Cells.Item(9,17).Value = 2.38
$ws.Cells.Item(9,18).Value = 1.17
$ws.Cells.Item(9,19).Value = 1.01
$ws.Cells.Item(9,20).Value = 1.01
$ws.Cells.Item(9,21).Value = 1.01
$ws.Cells.Item(9,22).Value = 1.46
$ws.Cells.Item(9,23).Value = 1.4
$ws.Cells.Item(9,24).Value = 9.2
$ws.Cells.Item(9,25).Value = 9.8
$ws.Cells.Item(9,26).Value = 16.5
$ws.Cells.Item(9,27).Value = 48
$ws.Cells.Item(9,28).Value = 10
$ws.Cells.Item(9,29).Value = 7
$ws.Cells.Item(9,30).Value = 14
$ws.Cells.Item(9,31).Value = 42
$ws.Cells.Item(9,32).Value = 19
$ws.Cells.Item(9,33).Value = 14.5
$ws.Cells.Item(9,34).Value = 24
$ws.Cells.Item(9,35).Value = 70
$ws.Cells.Item(9,36).Value = 60
$ws.Cells.Item(9,37).Value = 46
$ws.Cells.Item(9,38).Value = 75
$ws.Cells.Item(9,39).Value = 180
$ws.Cells.Item(9,40).Value = 60
$ws.Cells.Item(9,41).Value = 50

# Row 10: Brazilian Serie A | Fluminense vs Juventude
$ws.Cells.Item(10,1).Value = 'Brazilian Serie A'
$ws.Cells.Item(10,2).NumberFormat = '@'
$ws.Cells.Item(10,2).Value = '2025-10-16'
$ws.Cells.Item(10,2).Style = 'Normal'
$ws.Cells.Item(10,3).Value = '21:30:00'
$ws.Cells.Item(10,4).Value = 'Fluminense'
$ws.Cells.Item(10,5).Value = 'Juventude'
$ws.Cells.Item(10,6).Value = 1.37
$ws.Cells.Item(10,7).Value = 1.42
$ws.Cells.Item(10,8).Value = 3.4
$ws.Cells.Item(10,9).Value = 15
$ws.Cells.Item(10,10).Value = 5
$ws.Cells.Item(10,11).Value = 5.8
$ws.Cells.Item(10,12).Value = 1.44
$ws.Cells.Item(10,13).Value = 1.07
$ws.Cells.Item(10,14).Value = 1.25
$ws.Cells.Item(10,15).Value = 1.36
$ws.Cells.Item(10,16).Value = 1.25
$ws.Cells.Item(10,17).Value = 1.36
$ws.Cells.Item(10,18).Value = 1.19
$ws.Cells.Item(10,19).Value = 3.25
$ws.Cells.Item(10,20).Value = 1.01
$ws.Cells.Item(10,21).Value = 1.01
$ws.Cells.Item(10,22).Value = 1.07
$ws.Cells.Item(10,23).Value = 2.14
$ws.Cells.Item(10,24).Value = 14
$ws.Cells.Item(10,25).Value = 36
$ws.Cells.Item(10,26).Value = 130
$ws.Cells.Item(10,27).Value = 850
$ws.Cells.Item(10,28).Value = 8
$ws.Cells.Item(10,29).Value = 12.5
$ws.Cells.Item(10,30).Value = 55
$ws.Cells.Item(10,31).Value = 390
$ws.Cells.Item(10,32).Value = 6.8
$ws.Cells.Item(10,33).Value = 11.5
$ws.Cells.Item(10,34).Value = 46
$ws.Cells.Item(10,35).Value = 370
$ws.Cells.Item(10,36).Value = 10.5
$ws.Cells.Item(10,37).Value = 18.5
$ws.Cells.Item(10,38).Value = 65
$ws.Cells.Item(10,39).Value = 450
$ws.Cells.Item(10,40).Value = 10
$ws.Cells.Item(10,41).Value = 760

# Row 11: Honduras Liga Nacional | CD Olimpia vs CD Motagua
$ws.Cells.Item(11,1).Value = 'Honduras Liga Nacional'
$ws.Cells.Item(11,2).NumberFormat = '@'
$ws.Cells.Item(11,2).Value = '2025-10-16'
$ws.Cells.Item(11,2).Style = 'Normal'
$ws.Cells.Item(11,3).Value = '22:30:00'
$ws.Cells.Item(11,4).Value = 'CD Olimpia'
$ws.Cells.Item(11,5).Value = 'CD Motagua'
$ws.Cells.Item(11,6).Value = 1.04
$ws.Cells.Item(11,7).Value = 1000
$ws.Cells.Item(11,8).Value = 1.04
$ws.Cells.Item(11,9).Value = 1000
$ws.Cells.Item(11,10).Value = 1.02
$ws.Cells.Item(11,11).Value = 1000
$ws.Cells.Item(11,12).Value = 1.01
$ws.Cells.Item(11,13).Value = 1.01
$ws.Cells.Item(11,14).Value = 1.34
$ws.Cells.Item(11,15).Value = 1.28
$ws.Cells.Item(11,16).Value = 1.34
$ws.Cells.Item(11,17).Value = 1.28
$ws.Cells.Item(11,18).Value = 1.21
$ws.Cells.Item(11,19).Value = 1.28
$ws.Cells.Item(11,20).Value = 1.01
$ws.Cells.Item(11,21).Value = 1.01
$ws.Cells.Item(11,22).Value = 1.01
$ws.Cells.Item(11,23).Value = 1.01
$ws.Cells.Item(11,24).Value = 1000
$ws.Cells.Item(11,25).Value = 1000
$ws.Cells.Item(11,26).Value = 1000
$ws.Cells.Item(11,27).Value = 1000
$ws.Cells.Item(11,28).Value = 1000
$ws.Cells.Item(11,29).Value = 1000
$ws.Cells.Item(11,30).Value = 1000
$ws.Cells.Item(11,31).Value = 1000
$ws.Cells.Item(11,32).Value = 1000
$ws.Cells.Item(11,33).Value = 1000
$ws.Cells.Item(11,34).Value = 1000
$ws.Cells.Item(11,35).Value = 1000
$ws.Cells.Item(11,36).Value = 1000
$ws.Cells.Item(11,37).Value = 1000
$ws.Cells.Item(11,38).Value = 1000
$ws.Cells.Item(11,39).Value = 1000
$ws.Cells.Item(11,40).Value = 1000
$ws.Cells.Item(11,41).Value = 1000
